$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns remain plain text (not auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "59.050.57"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "2.506.29"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "537.51"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "138.40"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "0.561"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").Value = "2.502.03"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "5.41"
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  -4.11%  "
$ws.Range("D14").Value = "2.968.54"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "23.06"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").Value = "59.075.16"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "2.515.82"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").Value = "11.00"
$ws.Range("E19").Value = "  -2.91%  "
$ws.Range("D20").Value = "4.25"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "322.65"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").Value = "63.03"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").Value = "0.418"
$ws.Range("E25").Value = "  -5.10%  "
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("D28").Value = "7.69"
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("D29").Value = "0.0₃0771"
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("D30").Value = "6.65"
$ws.Range("E30").Value = "  -5.97%  "
$ws.Range("D31").Value = "1.78"
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").Value = "165.36"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "1.10"
$ws.Range("E34").Value = "  -11.61%  "
$ws.Range("D35").Value = "1.38"
$ws.Range("E35").Value = "  -8.15%  "
$ws.Range("D36").Value = "18.43"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("E37").Value = "  -7.68%  "
$ws.Range("D38").Value = "1.56"
$ws.Range("E38").Value = "  -4.72%  "
$ws.Range("D39").Value = "3.63"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("D40").Value = "0.801"
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("D41").Value = "5.20"
$ws.Range("E41").Value = "  -9.32%  "
$ws.Range("D42").Value = "278.94"
$ws.Range("E42").Value = "  -7.73%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "10.87"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.595"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "125.13"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "0.0936"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "0.0509"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").Value = "0.0221"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").Value = "17.62"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Value = "1.763.21"
$ws.Range("E51").Value = "  -3.12%  "

# Reset styling back to Normal so no visible/style side-effects remain from the text format
$ws.Range("D2:E51").Style = "Normal"
